$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eetu Pihamäki")

# New work-log entry for 9.11.2018 (row 28 of the "Eetu Pihamäki" sheet)
$ws.Range("A28").Value = (Get-Date -Year 2018 -Month 11 -Day 9 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("B28").Value = 0.41666666666666669
$ws.Range("C28").Value = 0.63194444444444442
$ws.Range("E28").Value = 3
$ws.Range("F28").Value = "4h Mm. Open ssl:llä salausavainten tekoa, itsekirjoitetun sertifikaatin luonti ja sen lisääminen keystore.jceks -tiedostoon. JVM-keystore jne. https://github.com/Eetu95/Open-source-IdM-solution/blob/master/Eetun%20muistiinpanoja/Ty%C3%B6t%20-%209.11.2018.txt"

# Row grew to fit the wrapped task description
$ws.Rows.Item(28).RowHeight = 90

# Move the active selection like the author left it
$ws.Activate()
$ws.Range("F29").Select()
